# Updates recomputed TPM-derived NATMI ligand-receptor metrics for Sema6b-Plxna2
# Replaces columns G:T (rows 2-26) with the refreshed values from the new TPM run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @(10.06074966666667, 30.182249, 0.5973429698540321, 0.607043146831649, 3, 1, 23.31637266666667, 69.949118, 0.5832682088864471, 0.650052653999578, 234.5801885340424, 2111.221696806382, 0.3484111641176723, 0.3946100086901689),
    @(10.06074966666667, 30.182249, 0.5973429698540321, 0.607043146831649, 3, 1, 3.113941333333333, 9.341823999999999, 0.07789646400133915, 0.08681564054027034, 31.32858423135288, 281.9572580821759, 0.04653090514768763, 0.05270083962777099),
    @(10.06074966666667, 30.182249, 0.5973429698540321, 0.607043146831649, 3, 1, 0.2259803333333333, 0.677941, 0.005652986686703996, 0.006300255941828, 2.273531563256555, 20.461784069309, 0.00337677185598107, 0.003824527192772063),
    @(10.06074966666667, 30.182249, 0.5973429698540321, 0.607043146831649, 2, 1, 12.3208515, 24.641703, 0.3082109335400439, 0.2290008064750633, 123.9570026216745, 743.742015730047, 0.1841076343822935, 0.1390133701896079),
    @(10.06074966666667, 30.182249, 0.5973429698540321, 0.607043146831649, 3.0, 1.0, 0.9982416666666666, 2.994725, 0.02497140688546588, 0.02783064304326019, 10.04305951516944, 90.38753563652499, 0.01491649435039761, 0.01689440113132901),
    @(2.694390666666667, 8.083172, 0.1599756853215596, 0.1625735102530455, 3, 1, 23.31637266666667, 69.949118, 0.5832682088864471, 0.650052653999578, 62.82341689358844, 565.4107520422959, 0.09330873144288795, 0.1056813418100199),
    @(2.694390666666667, 8.083172, 0.1599756853215596, 0.1625735102530455, 3, 1, 3.113941333333333, 9.341823999999999, 0.07789646400133915, 0.08681564054027034, 8.390174465080888, 75.51157018572799, 0.01246154021274043, 0.01411392342749836),
    @(2.694390666666667, 8.083172, 0.1599756853215596, 0.1625735102530455, 3, 1, 0.2259803333333333, 0.677941, 0.005652986686703996, 0.006300255941828, 0.6088793009835556, 5.479913708852, 0.0009043404193191242, 0.001024254723955585),
    @(2.694390666666667, 8.083172, 0.1599756853215596, 0.1625735102530455, 2, 1, 12.3208515, 24.641703, 0.3082109335400439, 0.2290008064750633, 33.197187286986, 199.183123721916, 0.04930625531666617, 0.03722946495942939),
    @(2.694390666666667, 8.083172, 0.1599756853215596, 0.1625735102530455, 3.0, 1.0, 0.9982416666666666, 2.994725, 0.02497140688546588, 0.02783064304326019, 2.689653029744444, 24.2068772677, 0.003994817929945916, 0.004524525332142312),
    @(1.673271333333333, 5.019813999999999, 0.09934815006247043, 0.1009614521127821, 3, 1, 23.31637266666667, 69.949118, 0.5832682088864471, 0.650052653999578, 39.01461798045022, 351.131561824052, 0.0579466175431191, 0.06563025989756534),
    @(1.673271333333333, 5.019813999999999, 0.09934815006247043, 0.1009614521127821, 3, 1, 3.113941333333333, 9.341823999999999, 0.07789646400133915, 0.08681564054027034, 5.210468766748443, 46.89421890073599, 0.007738869594940867, 0.008765033135047012),
    @(1.673271333333333, 5.019813999999999, 0.09934815006247043, 0.1009614521127821, 3, 1, 0.2259803333333333, 0.677941, 0.005652986686703996, 0.006300255941828, 0.3781264136637777, 3.403137722973999, 0.0005616137696518161, 0.0006360829885691388),
    @(1.673271333333333, 5.019813999999999, 0.09934815006247043, 0.1009614521127821, 2, 1, 12.3208515, 24.641703, 0.3082109335400439, 0.2290008064750633, 20.616127617207, 123.696765703242, 0.03062018607623038, 0.02312025395672059),
    @(1.673271333333333, 5.019813999999999, 0.09934815006247043, 0.1009614521127821, 3.0, 1.0, 0.9982416666666666, 2.994725, 0.02497140688546588, 0.02783064304326019, 1.670329164572222, 15.03296248115, 0.002480863078528271, 0.002809822134880047),
    @(0.8073985, 1.614797, 0.04793815906618067, 0.03247774718094421, 3, 1, 23.31637266666667, 69.949118, 0.5832682088864471, 0.650052653999578, 18.82560431650767, 112.953625899046, 0.0279608041758448, 0.0211122457509001),
    @(0.8073985, 1.614797, 0.04793815906618067, 0.03247774718094421, 3, 1, 3.113941333333333, 9.341823999999999, 0.07789646400133915, 0.08681564054027034, 2.514191561621333, 15.085149369728, 0.003734213081989212, 0.002819576424818631),
    @(0.8073985, 1.614797, 0.04793815906618067, 0.03247774718094421, 3, 1, 0.2259803333333333, 0.677941, 0.005652986686703996, 0.006300255941828, 0.1824561821628334, 1.094737092977, 0.0002709937749862178, 0.0002046181196539314),
    @(0.8073985, 1.614797, 0.04793815906618067, 0.03247774718094421, 2, 1, 12.3208515, 24.641703, 0.3082109335400439, 0.2290008064750633, 9.94783701982275, 39.791348079291, 0.01477506475797866, 0.007437430296929437),
    @(0.8073985, 1.614797, 0.04793815906618067, 0.03247774718094421, 3.0, 1.0, 0.9982416666666666, 2.994725, 0.02497140688546588, 0.02783064304326019, 0.8059788243041666, 4.835872945825, 0.001197083275381783, 0.0009038765886421083),
    @(1.606691, 4.820073, 0.09539503569575726, 0.09694414362157924, 3, 1, 23.31637266666667, 69.949118, 0.5832682088864471, 0.650052653999578, 37.46220611617933, 337.159855045614, 0.05564089160692303, 0.06301879785092385),
    @(1.606691, 4.820073, 0.09539503569575726, 0.09694414362157924, 3, 1, 3.113941333333333, 9.341823999999999, 0.07789646400133915, 0.08681564054027034, 5.003141514794666, 45.028273633152, 0.007430935963981019, 0.008416267925135366),
    @(1.606691, 4.820073, 0.09539503569575726, 0.09694414362157924, 3, 1, 0.2259803333333333, 0.677941, 0.005652986686703996, 0.006300255941828, 0.3630805677436666, 3.267725109693, 0.0005392668667657683, 0.0006107729168772817),
    @(1.606691, 4.820073, 0.09539503569575726, 0.09694414362157924, 2, 1, 12.3208515, 24.641703, 0.3082109335400439, 0.2290008064750633, 19.7958012173865, 118.774807304319, 0.02940179300687515, 0.02220028707237601),
    @(1.606691, 4.820073, 0.09539503569575726, 0.09694414362157924, 3.0, 1.0, 0.9982416666666666, 2.994725, 0.02497140688546588, 0.02783064304326019, 1.603865901658333, 14.434793114925, 0.002382148251212296, 0.002698017856266721)
)

$arr = New-Object "object[,]" $rows.Count, $rows[0].Count
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt $rows[$i].Count; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("G2:T26").Value = $arr

Write-Output "Updated $($rows.Count) rows x $($rows[0].Count) cols in range G2:T26"
